# Add "Physical Control" (4511.194 / M1 / Moving) to the Charges list.
# The data in Sheet1 is kept sorted alphabetically by "Name of Offense"
# (column A), so the new entry belongs on row 19, between
# "No Operator License - Never Held" (row 18) and
# "Possession of Marijuana Drug Paraphernalia" (the old row 19, which
# shifts down to row 20, along with every row after it).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Push rows 19..34 down one row, creating a blank row 19 for the new entry.
$ws.Rows.Item(19).Insert()

# Fill in the new charge on row 19.
$ws.Cells.Item(19, 1).Value = "Physical Control"
$ws.Cells.Item(19, 2).Value = "4511.194"
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 3).Value = "M1"
$ws.Cells.Item(19, 4).Value = "Moving"

# NOTE: the sheet is *mostly* sorted alphabetically by column A, but not
# perfectly (e.g. "Following Too Close" is already out of order in the
# source file) - so we must not re-run a full-column sort here, as that
# would reshuffle rows that the real edit left untouched. The new row is
# simply inserted directly at its correct alphabetical spot (row 19).

# Leave the selection on the newly added row, as in the authored edit.
$ws.Range("E19").Select()
